# The paragraph originally reads (visible text):
#   ... Le <m><pa>bled</pa> orbere</pa></m>  faict ...
# split across runs so that "b" and "led" are separate runs, then a
# Courier-New/blue run containing the literal text "</pa>", then a plain
# " " (space) run, then the "orbere" run.
#
# The target edit removes the "</pa>" run entirely (it was an erroneous
# stray closing tag) while keeping the single space that followed it, so
# the "led" run's text effectively becomes "led " and is immediately
# followed by "orbere".
#
# We locate the unique anchor text "led</pa> orbere" in the document,
# then compute the sub-range that covers just the "</pa>" substring
# inside that match and delete it. Word automatically merges the
# now-adjacent identically-formatted runs ("led" + " ") into a single
# "led " run, exactly matching the target OOXML.

$d = $word.ActiveDocument
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("led</pa> orbere", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $r.Start
    # Offsets within the matched text "led</pa> orbere":
    # l(0) e(1) d(2) <(3) /(4) p(5) a(6) >(7) ' '(8) o(9)...
    $tagStart = $matchStart + 3
    $tagEnd = $matchStart + 8
    $tagRange = $d.Range($tagStart, $tagEnd)
    $tagRange.Delete()
}
